# Merge the "<id>...<value>...</id>" three-run sequences (each split
# across a Courier-New "<id>" run, an Arial value run, and a Courier-New
# "</id>" run) into a single run per occurrence: "<id>value</id>".
#
# A Find/Replace where the replacement text equals the matched text
# forces the engine to rebuild the matched range as one run using the
# formatting of the first run in the match (Courier New), which is
# exactly the desired merged run.

$d = $word.ActiveDocument

$ids = @("p055r_4", "p055v_1", "p055v_2", "p055v_3", "p055v_4")

foreach ($idVal in $ids) {
    $needle = "<id>" + $idVal + "</id>"
    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
    # MatchAllWordForms, Forward, Wrap(1=wdFindContinue), Format,
    # ReplaceWith, Replace(2=wdReplaceAll)
    $d.Content.Find.Execute($needle, $true, $false, $false, $false, $false,
                             $true, 1, $false, $needle, 2)
}

